# KIOSC_Finance_Data.xlsx edit
# - Expenses!I4: "Invoiced" -> "Paid"
# - Expenses!L4: ""         -> "2025-04-24" (paymentDate, must stay text)
# - AuditLog: append two identical "UPDATE" audit rows (rows 4 & 5) recording
#   the status/paymentDate change on expense 7776af38-6226-4f6a-9d0e-68a9a692852f

$wb = $excel.ActiveWorkbook

# ---- Expenses sheet -------------------------------------------------
$expenses = $wb.Worksheets.Item("Expenses")

$expenses.Range("I4").Value = "Paid"

# The payment-date column would otherwise be auto-parsed into a date
# serial number by Excel's smart input; force it to stay plain text so it
# round-trips exactly like the source data ("2025-04-24").
$expenses.Range("L4").NumberFormat = "@"
$expenses.Range("L4").Value = "2025-04-24"

# ---- AuditLog sheet ---------------------------------------------------
$audit = $wb.Worksheets.Item("AuditLog")

$changeJson = '{"before":{"id":"7776af38-6226-4f6a-9d0e-68a9a692852f","date":"2023-02-15","description":"Educational Materials","supplier":"f99c571b-2d0d-48a8-b6fa-ec0a2bc491c9","amount":"3450","paymentType":"1","paymentCenter":"2","program":"2","status":"Invoiced","notes":"Materials for outreach program","invoiceDate":"2023-02-20","paymentDate":"","createdBy":"user","createdAt":"2023-02-15T14:00:00.000Z"},"after":{"id":"7776af38-6226-4f6a-9d0e-68a9a692852f","date":"2023-02-15","description":"Educational Materials","supplier":"f99c571b-2d0d-48a8-b6fa-ec0a2bc491c9","amount":"3450","paymentType":"1","paymentCenter":"2","program":"2","status":"Paid","notes":"Materials for outreach program","invoiceDate":"2023-02-20","paymentDate":"2025-04-24","createdBy":"user","createdAt":"2023-02-15T14:00:00.000Z"}}'

$auditId = "AUDIT1745477662005"
$entityType = "Expenses"
$entityId = "7776af38-6226-4f6a-9d0e-68a9a692852f"
$action = "UPDATE"
$userId = "1"
$username = "admin"
$timestamp = "2025-04-24T06:54:22.005Z"
$description = "Updated Expense 7776af38-6226-4f6a-9d0e-68a9a692852f"

foreach ($r in 4, 5) {
    $audit.Range("A$r").Value = $auditId
    $audit.Range("B$r").Value = $entityType
    $audit.Range("C$r").Value = $entityId
    $audit.Range("D$r").Value = $action

    # userId ("1") looks numeric and would otherwise be stored as a number;
    # the source file keeps every column as text, so force it.
    $audit.Range("E$r").NumberFormat = "@"
    $audit.Range("E$r").Value = $userId

    $audit.Range("F$r").Value = $username
    $audit.Range("G$r").Value = $timestamp
    $audit.Range("H$r").Value = $changeJson
    $audit.Range("I$r").Value = $description
}
